$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.793.27'
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.635.28'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.97'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("E8").Value = '  +0.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06421'
$ws.Range("E9").Value = '  +1.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.24'
$ws.Range("E10").Value = '  +3.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07778'
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.862.25'
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.633.64'
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5585'
$ws.Range("E15").Value = '  +1.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅7649'
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.24'
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.807.18'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.370'
$ws.Range("E20").Value = '  -1.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.80'
$ws.Range("E21").Value = '  -1.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.904'
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.137'
$ws.Range("E23").Value = '  +1.82%  '
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.775'
$ws.Range("E25").Value = '  -6.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '138.92'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1229'
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.812'
$ws.Range("E28").Value = '  +0.74%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04935'
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.292'
$ws.Range("E32").Value = '  +1.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.254'
$ws.Range("E33").Value = '  +2.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.570'
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9020'
$ws.Range("E36").Value = '  +0.65%  '
$ws.Range("E37").Value = '  +0.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5561'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.130.75'
$ws.Range("E39").Value = '  +1.31%  '
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9953'
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.449'
$ws.Range("E42").Value = '  -2.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '98.83'
$ws.Range("E43").Value = '  +1.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7987'
$ws.Range("E44").Value = '  +0.45%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0₈115'
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.58'
$ws.Range("E46").Value = '  +1.46%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4256'
$ws.Range("E47").Value = '  -4.09%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.771'
$ws.Range("E48").Value = '  +2.48%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05030'
$ws.Range("E49").Value = '  -2.07%  '
$ws.Range("B50").Value = 'Frax'
$ws.Range("C50").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9969'
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.002'
$ws.Range("E51").Value = '  +0.17%  '
